$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E3").Value = 15.9729
$ws.Range("A4").Value = -20.85499999999998
$ws.Range("A6").Value = -22.71130000000002
$ws.Range("A7").Value = -21.96010000000001
$ws.Range("B7").Value = 4.799900000000003
$ws.Range("A8").Value = -22.30080000000002
$ws.Range("B11").Value = 5.410599999999999
$ws.Range("B12").Value = 5.1019
$ws.Range("D12").Value = -6.044599999999997
$ws.Range("E12").Value = 18.75910000000001
$ws.Range("D13").Value = -8.451000000000002
$ws.Range("E13").Value = 16.55690000000001
$ws.Range("D14").Value = -8.2387
$ws.Range("B15").Value = 5.036099999999998
$ws.Range("A16").Value = -21.48279999999999
$ws.Range("D16").Value = -9.128400000000006
$ws.Range("D19").Value = -7.531899999999991
$ws.Range("A20").Value = -22.51260000000001
$ws.Range("B20").Value = 4.366299999999996
$ws.Range("D20").Value = -7.665199999999995
$ws.Range("A21").Value = -22.2062
$ws.Range("B21").Value = 5.383799999999997
$ws.Range("B22").Value = 9.940600000000002
$ws.Range("D22").Value = -6.825299999999991
$ws.Range("E22").Value = 16.99910000000001
$ws.Range("B23").Value = 8.895200000000006
$ws.Range("E25").Value = 17.0746
$ws.Range("A28").Value = -22.35340000000001
$ws.Range("A29").Value = -21.7276
$ws.Range("B29").Value = 5.388500000000001
$ws.Range("E29").Value = 17.29140000000001
$ws.Range("A30").Value = -21.88450000000001
$ws.Range("A32").Value = -21.20699999999998
$ws.Range("B34").Value = 9.727400000000008
$ws.Range("E34").Value = 17.1845
$ws.Range("D36").Value = -8.361299999999998
$ws.Range("A40").Value = -19.18239999999999
$ws.Range("B42").Value = 9.851299999999997
$ws.Range("B43").Value = 5.078900000000004
$ws.Range("D43").Value = -8.443900000000001
$ws.Range("E43").Value = 16.55049999999999
$ws.Range("B44").Value = 5.489899999999997
$ws.Range("B45").Value = 5.017000000000003
$ws.Range("A46").Value = -22.0647
$ws.Range("B46").Value = 5.7593
$ws.Range("D46").Value = -8.775299999999994
$ws.Range("E48").Value = 17.45890000000001
$ws.Range("B50").Value = 4.556799999999995
$ws.Range("D50").Value = -8.195799999999997
$ws.Range("A51").Value = -22.14879999999999
$ws.Range("B51").Value = 5.317199999999998
$ws.Range("A52").Value = -22.13939999999999
$ws.Range("A57").Value = -22.71610000000002
$ws.Range("B57").Value = 5.205699999999994
$ws.Range("A59").Value = -22.4417
$ws.Range("E60").Value = 15.38250000000001
$ws.Range("A62").Value = -22.23010000000001
$ws.Range("B65").Value = 5.428299999999999
$ws.Range("A66").Value = -21.52590000000001
$ws.Range("B66").Value = 4.885499999999997
$ws.Range("B67").Value = 5.299400000000003
$ws.Range("E68").Value = 17.45330000000001
$ws.Range("E70").Value = 18.14990000000002
$ws.Range("E71").Value = 17.03260000000001
$ws.Range("A73").Value = -19.86799999999998
$ws.Range("E73").Value = 17.28470000000002
$ws.Range("A74").Value = -21.97049999999998
$ws.Range("D76").Value = -7.988700000000001
$ws.Range("A77").Value = -19.80349999999999
$ws.Range("E78").Value = 17.07980000000001
$ws.Range("B79").Value = 9.737000000000007
$ws.Range("B84").Value = 5.498000000000001
$ws.Range("B87").Value = 5.448299999999997
$ws.Range("E87").Value = 16.3571
$ws.Range("A92").Value = -21.35270000000002
$ws.Range("B92").Value = 5.269299999999998
$ws.Range("E92").Value = 18.97560000000002
$ws.Range("D95").Value = -8.143700000000006
$ws.Range("B97").Value = 6.382499999999997
$ws.Range("D97").Value = -8.380299999999998
$ws.Range("D99").Value = -8.289300000000001
$ws.Range("A100").Value = -22.2128
$ws.Range("E101").Value = 17.07650000000001
